# Regenerate s_vals data to filter save games.
# Updates the numeric B:G columns (TB, d2S, K, IP, Win, sum) for rows 2-9
# on the active worksheet, leaving the A column dates and F (Win) column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (B, C, D, E, G) new values. F (Win) is unchanged by this edit.
$data = @{
    2 = @(0.127881588408715,   0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0.8245905782990299)
    3 = @(3.230985683306322,   1.667794583268128,   0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(1.459612070389937,   1.667794583268128,   0.1575252929769615, 0.496779210170732, 3.781711156805759)
    5 = @(3.230985683306322,   1.667794583268128,   0.8054896365839992, 0.496779210170732, 6.201049113329182)
    6 = @(0.127881588408715,   0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.472554921906068)
    7 = @(1.459612070389937,   1.667794583268128,   0.8054896365839992, 0.496779210170732, 4.429675500412797)
    8 = @(3.230985683306322,   1.667794583268128,   0.8054896365839992, 0.496779210170732, 6.201049113329182)
    9 = @(3.230985683306322,   1.667794583268128,   0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G - sum
}
